$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.075.92"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "1.833.34"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "324.31"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4642"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "0.3871"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "0.07851"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'0.9610"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "1.816.08"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "5.691"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "6.909"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "0.06847"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "88.48"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.000009934"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "16.66"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "28.073.52"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").Value = "5.307"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").Value = "2.094"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "2.036.85"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "154.67"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").Value = "19.14"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'5.660"
$ws.Range("E28").Value = "  -6.41%  "
$ws.Range("D29").Value = "1.957"
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").Value = "'118.30"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "0.9353"
$ws.Range("E31").Value = "  -4.27%  "
$ws.Range("D32").Value = "0.09235"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "5.256"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").Value = "'1.320"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").Value = "3.306"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("D36").Value = "'0.05860"
$ws.Range("E36").Value = "  -4.88%  "
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").Value = "1.148"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").Value = "7.751"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").Value = "0.5593"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").Value = "9.888"
$ws.Range("D42").Value = "0.1761"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").Value = "0.07244"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").Value = "'11.60"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "0.5267"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").Value = "1.146"
$ws.Range("E46").Value = "  -8.25%  "
$ws.Range("D47").Value = "2.101"
$ws.Range("E47").Value = "  -11.52%  "
$ws.Range("D48").Value = "1.822"
$ws.Range("E48").Value = "  -4.75%  "
$ws.Range("D49").Value = "112.13"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "1.025"
$ws.Range("E51").Value = "  +0.33%  "
